$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.936.46'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.764.52'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.81'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4681'
$ws.Range('E7').Value = '  +1.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3519'
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.20'
$ws.Range('E9').Value = '  +3.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07382'
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.083'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.62'
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.002'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.174'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.758.98'
$ws.Range('E16').Value = '  -0.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.34'
$ws.Range('E17').Value = '  -1.39%  '
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06421'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.90'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.965.60'
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.151'
$ws.Range('E25').Value = '  +3.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.24'
$ws.Range('E26').Value = '  -1.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.00'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.966.31'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.179'
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.87'
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.072'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09293'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.645'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.551'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.69'
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02267'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06078'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2063'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.906'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6143'
$ws.Range('E40').Value = '  -2.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.187'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.401'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.753'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.10'
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.739'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5788'
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '123.20'
$ws.Range('E47').Value = '  +0.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.931'
$ws.Range('E48').Value = '  -0.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06821'
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.14'
$ws.Range('E51').Value = '  +0.10%  '
